# Refreshed KPI_Summary metrics (relative-path data pull produced new pre/post
# test-vs-control aggregates) for all KPI rows on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: LPE
$ws.Range("B2").Value = 0.2547949735449735
$ws.Range("C2").Value = 0.263957871191955
$ws.Range("D2").Value = 0.2132826278659612
$ws.Range("E2").Value = 0.2255245715435933
$ws.Range("F2").Value = -0.04151234567901235
$ws.Range("G2").Value = -0.03843329964836176
$ws.Range("H2").Value = -0.003079046030650581
$ws.Range("I2").Value = -16.29245078953062
$ws.Range("J2").Value = -14.56039157870472
$ws.Range("K2").Value = -1.732059210825902
$ws.Range("L2").Value = "↓ Worse"
$ws.Range("M2").Value = -0.003079046030650706
$ws.Range("N2").Value = -0.08089603611640606
$ws.Range("O2").Value = 0.9366691722751022

# Row 3: avg_loan_size
$ws.Range("B3").Value = 7147.122709876543
$ws.Range("C3").Value = 7275.474526538211
$ws.Range("D3").Value = 7264.157952601411
$ws.Range("E3").Value = 6833.485441109632
$ws.Range("F3").Value = 117.0352427248675
$ws.Range("G3").Value = -441.9890854285795
$ws.Range("H3").Value = 559.024328153447
$ws.Range("I3").Value = 1.637515507648108
$ws.Range("J3").Value = -6.075055088384513
$ws.Range("K3").Value = 7.712570596032622
$ws.Range("M3").Value = 559.0243281534467
$ws.Range("N3").Value = 0.4552593390047078
$ws.Range("O3").Value = 0.6593739956813148

# Row 4: dq29_pot30_payment_rate_$_up_to_day
$ws.Range("B4").Value = 0.0006399292330040934
$ws.Range("C4").Value = -0.000004846238655923245
$ws.Range("D4").Value = 0.004293299205260926
$ws.Range("E4").Value = 0.002000721913352637
$ws.Range("F4").Value = 0.003653369972256832
$ws.Range("G4").Value = 0.002005568152008561
$ws.Range("H4").Value = 0.001647801820248271
$ws.Range("I4").Value = 570.9021847785259
$ws.Range("J4").Value = 41384.01540661404
$ws.Range("K4").Value = -40813.11322183551
$ws.Range("M4").Value = 0.001647801820248272
$ws.Range("N4").Value = 1.242862057599441
$ws.Range("O4").Value = 0.2434738379017352

# Row 5: dq29_pot30_payment_rate_unit_per_day
$ws.Range("B5").Value = 0.006899202885103575
$ws.Range("C5").Value = 0.007834284629202642
$ws.Range("D5").Value = 0.008360770404051872
$ws.Range("E5").Value = 0.007925903806162591
$ws.Range("F5").Value = 0.001461567518948298
$ws.Range("G5").Value = 0.00009161917695994894
$ws.Range("H5").Value = 0.001369948341988349
$ws.Range("I5").Value = 21.18458528164236
$ws.Range("J5").Value = 1.169464492245214
$ws.Range("K5").Value = 20.01512078939715
$ws.Range("M5").Value = 0.001369948341988348
$ws.Range("N5").Value = 1.047490883728273
$ws.Range("O5").Value = 0.3187924072906324

# Row 6: dq29_pot30_payment_rate_unit_up_to_day
$ws.Range("B6").Value = 0.7671867262947433
$ws.Range("C6").Value = 0.7105722547456687
$ws.Range("D6").Value = 0.4601640860580402
$ws.Range("E6").Value = 0.4818816658667086
$ws.Range("F6").Value = -0.3070226402367031
$ws.Range("G6").Value = -0.2286905888789601
$ws.Range("H6").Value = -0.07833205135774304
$ws.Range("I6").Value = -40.01928470784685
$ws.Range("J6").Value = -32.18400202817011
$ws.Range("K6").Value = -7.835282679676737
$ws.Range("M6").Value = -0.07833205135774296
$ws.Range("N6").Value = -3.961831484196368
$ws.Range("O6").Value = 0.003435097900523567

# Row 7: dq30_pct_$
$ws.Range("B7").Value = 0.9989945035151711
$ws.Range("C7").Value = 0.9985603325721083
$ws.Range("D7").Value = 0.9960883203883735
$ws.Range("E7").Value = 0.9978918602057362
$ws.Range("F7").Value = -0.002906183126797357
$ws.Range("G7").Value = -0.0006684723663720845
$ws.Range("H7").Value = -0.002237710760425273
$ws.Range("I7").Value = -0.2909108224891695
$ws.Range("J7").Value = -0.06694361317660878
$ws.Range("K7").Value = -0.2239672093125607
$ws.Range("M7").Value = -0.002237710760425449
$ws.Range("N7").Value = -1.74849548988991
$ws.Range("O7").Value = 0.1037774793690867

# Row 8: dq30_pct_unit
$ws.Range("B8").Value = 0.05265220743094221
$ws.Range("C8").Value = 0.04897191210739486
$ws.Range("D8").Value = 0.03306492644480295
$ws.Range("E8").Value = 0.03270223026851257
$ws.Range("F8").Value = -0.01958728098613925
$ws.Range("G8").Value = -0.01626968183888229
$ws.Range("H8").Value = -0.003317599147256954
$ws.Range("I8").Value = -37.20125317030557
$ws.Range("J8").Value = -33.22247618839766
$ws.Range("K8").Value = -3.978776981907906
$ws.Range("M8").Value = -0.003317599147256961
$ws.Range("N8").Value = -1.329733367773564
$ws.Range("O8").Value = 0.2116193810713998
